$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 18 de Julio de 2020 a las 08:53"

# Row 6 - India
$ws.Range("B6").Value = 1040746
$ws.Range("C6").Value = 289
$ws.Range("D6").Value = 654130
$ws.Range("E6").Value = 360325
$ws.Range("G6").Value = 6
$ws.Range("H6").Value = 26291

# Row 51 - Afganistan
$ws.Range("B51").Value = 35289
$ws.Range("C51").Value = 60
$ws.Range("D51").Value = 23243
$ws.Range("E51").Value = 10899

# Row 52 - Armenia
$ws.Range("B52").Value = 34462
$ws.Range("C52").Value = 461
$ws.Range("D52").Value = 23123
$ws.Range("E52").Value = 10708
$ws.Range("G52").Value = 11
$ws.Range("H52").Value = 631

# Row 75 - El Salvador
$ws.Range("D75").Value = 6441
$ws.Range("E75").Value = 4442
$ws.Range("G75").Value = 15
$ws.Range("H75").Value = 324

# Row 146 - Georgia
$ws.Range("B146").Value = 1018
$ws.Range("C146").Value = 8
$ws.Range("D146").Value = 895
$ws.Range("E146").Value = 108

# Row 191 - Islas Turcas y Caicos
$ws.Range("B191").Value = 75
$ws.Range("C191").Value = 1
$ws.Range("E191").Value = 61

# Row 209 - Papua Nueva Guinea
$ws.Range("B209").Value = 16
$ws.Range("E209").Value = 8
